$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: replace the single inline-drawing paragraph at $ParaIndex by pulling
# its WordOpenXML, editing the <w:p>...</w:p> fragment with plain string
# substitutions, and feeding the fixed-up fragment back in via InsertXML
# (which cleanly replaces the paragraph's contents when the fragment is
# exactly one <w:p>...</w:p> element).
# ---------------------------------------------------------------------------
function Update-ImageParagraph {
    param(
        [int]$ParaIndex,
        [string]$Marker,
        [string]$EmbedTarget,     # the correct final r:embed value, e.g. "rId6"
        [string[]]$Replacements   # pairs: old1, new1, old2, new2, ...
    )

    $p = $d.Paragraphs.Item($ParaIndex)
    $rng = $p.Range
    $xmlFull = $rng.WordOpenXML

    $markerPos = $xmlFull.IndexOf($Marker)
    $before = $xmlFull.Substring(0, $markerPos)
    $pStart = $before.LastIndexOf("<w:p ")

    $after = $xmlFull.Substring($markerPos)
    $endRel = $after.IndexOf("</w:p>")
    $pEnd = $markerPos + $endRel + 6

    $pXml = $xmlFull.Substring($pStart, $pEnd - $pStart)

    for ($i = 0; $i -lt $Replacements.Length; $i += 2) {
        $pXml = $pXml.Replace($Replacements[$i], $Replacements[$i + 1])
    }

    # Range.WordOpenXML re-derives a fresh, snippet-local relationship-id
    # numbering, so whatever r:embed value comes back is not necessarily the
    # real document.xml.rels id for this image. Force it to the known-correct
    # value (regardless of what WordOpenXML happened to assign) via regex,
    # since each paragraph fragment contains exactly one blip.
    $pXml = [regex]::Replace($pXml, 'r:embed="[^"]*"', 'r:embed="' + $EmbedTarget + '"')

    $rng.InsertXML($pXml)
}

# NOTE: this shim's PowerShell only binds *positional* arguments correctly;
# named (-Param value) binding silently drops the value. So every call below
# is positional: Update-ImageParagraph <ParaIndex> <Marker> <Replacements[]>

# Each call's 3rd positional arg is the *true* document.xml.rels id that
# this image's <a:blip r:embed="..."/> must end up with - Update-ImageParagraph
# forces it there directly, independent of whatever id WordOpenXML reports.

# --- Block 1 : image1.png -> image2.png (docPr/cNvPr rename only) ----------
Update-ImageParagraph 9 "image1.png" "rId6" @(
    '<wp:docPr id="3" name="image1.png"/>', '<wp:docPr id="4" name="image2.png"/>',
    '<pic:cNvPr id="0" name="image1.png"/>', '<pic:cNvPr id="0" name="image2.png"/>'
)

# --- Block 2 : image6.png -> image4.png, resize 1238250 -> 1435100 ---------
Update-ImageParagraph 16 "image6.png" "rId7" @(
    '<wp:extent cx="1238250" cy="1238250"/>', '<wp:extent cx="1435100" cy="1435100"/>',
    '<wp:docPr id="6" name="image6.png"/>', '<wp:docPr id="1" name="image4.png"/>',
    '<pic:cNvPr id="0" name="image6.png"/>', '<pic:cNvPr id="0" name="image4.png"/>',
    '<a:ext cx="1238250" cy="1238250"/>', '<a:ext cx="1435100" cy="1435100"/>'
)

# --- Block 3 : image4.png -> image1.png (docPr/cNvPr rename only) ----------
Update-ImageParagraph 24 "image4.png" "rId8" @(
    '<wp:docPr id="5" name="image4.png"/>', '<wp:docPr id="6" name="image1.png"/>',
    '<pic:cNvPr id="0" name="image4.png"/>', '<pic:cNvPr id="0" name="image1.png"/>'
)

# --- Block 4 : image7.jpg -> image5.jpg (docPr/cNvPr rename only) ----------
Update-ImageParagraph 33 "image7.jpg" "rId9" @(
    '<wp:docPr id="4" name="image7.jpg"/>', '<wp:docPr id="5" name="image5.jpg"/>',
    '<pic:cNvPr id="0" name="image7.jpg"/>', '<pic:cNvPr id="0" name="image5.jpg"/>'
)

# --- Block 5 : image3.jpg -> image6.jpg (docPr/cNvPr rename only) ----------
Update-ImageParagraph 43 "image3.jpg" "rId10" @(
    '<wp:docPr id="2" name="image3.jpg"/>', '<wp:docPr id="3" name="image6.jpg"/>',
    '<pic:cNvPr id="0" name="image3.jpg"/>', '<pic:cNvPr id="0" name="image6.jpg"/>'
)

# --- Block 6 : image2.png -> image4.png, blip rId11 -> rId7 -----------------
Update-ImageParagraph 46 "image2.png" "rId7" @(
    '<wp:docPr id="1" name="image2.png"/>', '<wp:docPr id="2" name="image4.png"/>',
    '<pic:cNvPr id="0" name="image2.png"/>', '<pic:cNvPr id="0" name="image4.png"/>'
)

# --- Block 7 : image5.png -> image3.png, blip rId12 -> rId11 ----------------
Update-ImageParagraph 52 "image5.png" "rId11" @(
    '<wp:docPr id="7" name="image5.png"/>', '<wp:docPr id="7" name="image3.png"/>',
    '<pic:cNvPr id="0" name="image5.png"/>', '<pic:cNvPr id="0" name="image3.png"/>'
)

# --- Speaker name text change -----------------------------------------------
$d.Content.Find.Execute("Susana López ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Angel Moreno", 2)
